$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.151.99'
$ws.Range("E2").Value = '  -0.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.136.68'
$ws.Range("E3").Value = '  -1.41%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.16'
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.51'
$ws.Range("E6").Value = '  -3.89%  '

$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.572'
$ws.Range("E8").Value = '  -5.94%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.149.94'
$ws.Range("E9").Value = '  -0.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.117'
$ws.Range("E10").Value = '  -3.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.57'
$ws.Range("E11").Value = '  -3.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.384'
$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.689.03'
$ws.Range("E13").Value = '  -1.21%  '

$ws.Range("E14").Value = '  -0.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.225.64'
$ws.Range("E15").Value = '  -0.51%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.05'
$ws.Range("E16").Value = '  -1.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.143.61'
$ws.Range("E17").Value = '  -1.00%  '

$ws.Range("E18").Value = '  -3.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '399.99'
$ws.Range("E19").Value = '  -4.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.22'
$ws.Range("E20").Value = '  -2.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.49'
$ws.Range("E21").Value = '  -3.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.07'
$ws.Range("E22").Value = '  -0.64%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.02'
$ws.Range("E24").Value = '  -2.98%  '

$ws.Range("E25").Value = '  -1.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.195'
$ws.Range("E26").Value = '  -5.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000100'
$ws.Range("E27").Value = '  -5.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.76'
$ws.Range("E28").Value = '  -1.31%  '

$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.80'
$ws.Range("E31").Value = '  -2.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.11'
$ws.Range("E32").Value = '  -2.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '160.32'
$ws.Range("E33").Value = '  +1.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.25'
$ws.Range("E34").Value = '  -1.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.80'
$ws.Range("E35").Value = '  -4.96%  '

$ws.Range("E36").Value = '  -3.51%  '

$ws.Range("E37").Value = '  -2.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.648.35'
$ws.Range("E38").Value = '  -2.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.67'
$ws.Range("E39").Value = '  -2.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.60'
$ws.Range("E40").Value = '  -3.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.05'
$ws.Range("E41").Value = '  -3.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.29'
$ws.Range("E42").Value = '  -2.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.690'
$ws.Range("E43").Value = '  -2.92%  '

$ws.Range("E44").Value = '  -2.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.44'
$ws.Range("E45").Value = '  -3.64%  '

$ws.Range("E46").Value = '  -3.89%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.97'
$ws.Range("E47").Value = '  -3.16%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '284.57'
$ws.Range("E48").Value = '  -3.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.997'
$ws.Range("E49").Value = '  -0.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0972'
$ws.Range("E50").Value = '  -1.69%  '

$ws.Range("E51").Value = '  +0.20%  '
